$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (71) with the latest "cotação" (quote) values, following
# the same layout/style as the preceding rows (date in column A formatted
# like the other date cells, values B:E as text strings with comma decimals).

$newRow = 71

$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = 45975
$dateCell.NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = "22,2279"
$ws.Cells.Item($newRow, 3).Value = "15,8858"
$ws.Cells.Item($newRow, 4).Value = "15,6322"
$ws.Cells.Item($newRow, 5).Value = "15,6322"
